$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Month" column values from "May" to "June" (all data rows, col B, rows 2-46)
$ws.Range("B2:B46").Value = "June"

# Update "Bioassay" column values from 1 to 2 (col A, rows 2-46)
$ws.Range("A2:A46").Value = 2

# Swap the FvFm values between rows 8 and 9 (column F)
$f8 = $ws.Range("F8").Value2
$f9 = $ws.Range("F9").Value2
$ws.Range("F8").Value = $f9
$ws.Range("F9").Value = $f8

# Update the active selection to F11
$ws.Range("F11").Select()
